# "body2 and 3 最终版"
#
# The document had two near-duplicate paragraphs about the economic
# reasons for living with aging parents: an earlier, struck-through
# draft and a revised, clean version directly after it. This edit
# removes the struck-through draft paragraph entirely (it merges into
# the paragraph that follows, since deleting the paragraph mark joins
# the two) and then tightens the wording of the surviving paragraph:
#   - "It is cheaper for two families to live in one home than to have
#      a home respectively." -> "For two families, it is cheaper to
#      live in one home than to have a home respectively."
#   - "Caring for them at home seems more convenient and cheaper than
#      paying ..." -> "caring for them at home seems more convenient
#      and economical than paying ..." (lower-cased "caring" since it
#      no longer starts a sentence, and "cheaper" -> "economical").

$d = $word.ActiveDocument

# --- 1. Delete the whole struck-through draft paragraph -------------
# It begins "In addition to traditional culture, economy ..." and ends
# "... or to consider a move into assisted living." Locate it via the
# unique "for each to have a separate home" phrase (only present in
# this paragraph) and delete the entire paragraph range, including its
# end-of-paragraph mark, so the text that used to follow it becomes
# the start of a single merged paragraph.
$targetIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*for each to have a separate home*") {
        $targetIndex = $i
    }
}
if ($targetIndex -gt 0) {
    $para = $d.Paragraphs.Item($targetIndex)
    $paraRange = $d.Range($para.Range.Start, $para.Range.End)
    $paraRange.Delete()
}

# --- 2. Rework "It is cheaper for two families ..." ------------------
$d.Content.Find.Execute(
    "It is cheaper for two families to live in one home than to have a home",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "For two families, it is cheaper to live in one home than to have a home",
    2)

# --- 3. Rework "Caring for them at home seems more convenient and ----
#        cheaper than paying ..."
$d.Content.Find.Execute(
    "our parents. Furthermore, Caring for them at home seems more convenient and cheaper than paying",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "our parents. Furthermore, caring for them at home seems more convenient and economical than paying",
    2)
